$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = $ws1.Range("A1").Value2
$text = $text -replace [regex]::Escape("1000 Bs = 9.17 = 38176.15 pesos"), "1000 Bs = 9.12 = 37966.71 pesos"
$text = $text -replace [regex]::Escape("38176.15 pesos = 9.15 = 972.26 Bs"), "37966.71 pesos = 9.07 = 963.69 Bs"
$ws1.Range("A1").Value = $text

# --- Sheet "tasas": update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 109.65
$ws2.Range("O10").Value = 4163.05
$ws2.Range("N12").Value = 4183.98
